$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "29.239.96"
$ws.Range("E2").Value = "  +2.62%  "

$ws.Range("D3").Value = "1.891.33"
$ws.Range("E3").Value = "  +0.58%  "

$ws.Range("E4").Value = "  -0.83%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.00"
$ws.Range("E5").Value = "  -0.52%  "

$ws.Range("E6").Value = "  -0.80%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5137"
$ws.Range("E7").Value = "  +0.42%  "

$ws.Range("E8").Value = "  -1.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08407"
$ws.Range("E9").Value = "  -0.21%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.38"
$ws.Range("E10").Value = "  +1.51%  "

$ws.Range("E11").Value = "  +0.15%  "

$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.222"
$ws.Range("E12").Value = "  -0.81%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.889.32"
$ws.Range("E13").Value = "  +0.53%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.60"
$ws.Range("E14").Value = "  +0.40%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.297"
$ws.Range("E15").Value = "  +0.18%  "

$ws.Range("E16").Value = "  -0.82%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.95"
$ws.Range("E17").Value = "  +1.85%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001103"
$ws.Range("E18").Value = "  -0.54%  "

$ws.Range("E19").Value = "  -0.26%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.78"
$ws.Range("E20").Value = "  +0.35%  "

$ws.Range("E21").Value = "  -0.86%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.996"
$ws.Range("E22").Value = "  +0.60%  "

$ws.Range("D23").Value = "29.251.06"
$ws.Range("E23").Value = "  +2.48%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.11"
$ws.Range("E24").Value = "  -0.25%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.213"
$ws.Range("E25").Value = "  -2.56%  "

$ws.Range("D26").Value = "2.109.08"
$ws.Range("E26").Value = "  +0.63%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.10"
$ws.Range("E27").Value = "  -1.35%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.91"
$ws.Range("E28").Value = "  +0.70%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.425"
$ws.Range("E29").Value = "  +1.58%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.72"
$ws.Range("E30").Value = "  +0.15%  "

$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1045"
$ws.Range("E31").Value = "  -0.62%  "

$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.055"
$ws.Range("E32").Value = "  +0.44%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.136"
$ws.Range("E33").Value = "  +5.93%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.656"
$ws.Range("E34").Value = "  +1.18%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02473"
$ws.Range("E35").Value = "  +1.25%  "

$ws.Range("E36").Value = "  +0.41%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.028"
$ws.Range("E37").Value = "  +1.03%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2184"
$ws.Range("E38").Value = "  -0.18%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.222"
$ws.Range("E39").Value = "  +2.55%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.128"
$ws.Range("E40").Value = "  +0.75%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6488"
$ws.Range("E41").Value = "  +0.47%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.228"
$ws.Range("E42").Value = "  -3.45%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.23"
$ws.Range("E43").Value = "  +0.31%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6038"
$ws.Range("E44").Value = "  -0.63%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.17"
$ws.Range("E45").Value = "  +0.91%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.678"
$ws.Range("E46").Value = "  -0.97%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.039"
$ws.Range("E47").Value = "  +1.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.227"
$ws.Range("E48").Value = "  +1.33%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "122.85"
$ws.Range("E49").Value = "  +0.31%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.177"
$ws.Range("E50").Value = "  -2.24%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "77.40"
$ws.Range("E51").Value = "  +0.29%  "
